# Daily attendance processing - 2025-11-14 23:21:31
#
# The "Recorded By" column (G) lists the users/processes that touched each
# attendance row, separated by ", ". Re-processing promotes the actual
# reporting users (backup@backdoor.com / dnasr281@gmail.com) ahead of the
# generic "System"/"system" (and other) entries, while keeping each group's
# relative order stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Users that should be sorted to the front of the "Recorded By" list.
$priority = @{}
$priority["backup@backdoor.com"] = 0
$priority["dnasr281@gmail.com"] = 0

# Find the last populated row in column A (mirrors the sheet's used range).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -eq $val) { continue }
    if ($val.GetType().Name -ne "String") { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    $front = @()
    $back = @()
    foreach ($p in $parts) {
        if ($priority.ContainsKey($p)) {
            $front += $p
        } else {
            $back += $p
        }
    }

    $reordered = $front + $back
    $joined = $reordered -join ", "

    if ($joined -ne $val) {
        $cell.Value = $joined
    }
}
